$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 390; existing rows 390-428 shift down to 391-429
$ws.Rows.Item(390).Insert()

# Populate the newly inserted row 390 with the new weekly record
$ws.Cells.Item(390, 1).Value = 8
$ws.Cells.Item(390, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(390, 3).Value = "Coquimbo"
$ws.Cells.Item(390, 4).Value = (Get-Date -Year 2023 -Month 8 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(390, 5).Value = 4
$ws.Cells.Item(390, 6).Value = 100112031
$ws.Cells.Item(390, 7).Value = "Poroto verde"
$ws.Cells.Item(390, 8).Value = "Magnum"
$ws.Cells.Item(390, 9).Value = "Primera"
$ws.Cells.Item(390, 10).Value = 360
$ws.Cells.Item(390, 11).Value = 28000
$ws.Cells.Item(390, 12).Value = 29000
$ws.Cells.Item(390, 13).Value = 28500
$ws.Cells.Item(390, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(390, 15).Value = "Perú"
$ws.Cells.Item(390, 16).Value = 1140
$ws.Cells.Item(390, 17).Value = 25
$ws.Cells.Item(390, 18).Value = "Hortaliza"
